$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Based off their project plan and using the software chosen, students are to present annotated screenshots showing the following:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Based off their project plan and using the software chosen, students are to present annotated screenshots showing the following:",
    2)
